$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Milestone 1 -> "M" with start/end dates
$ws.Range("C5").Value = "M"
$ws.Range("D5").Value = "01/01"
$ws.Range("E5").Value = "01/14"

# Sub-activities referencing milestone 1
$ws.Range("C6").Value = "M - A"
$ws.Range("C7").Value = "M - B"

# Milestone 2 -> "N" with start/end dates
$ws.Range("C8").Value = "N"
$ws.Range("D8").Value = "01/15"
$ws.Range("E8").Value = "01/28"

# Sub-activities referencing milestone 2
$ws.Range("C9").Value = "N - C"
$ws.Range("C10").Value = "N - D"
